$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stable donor cells for text placeholders (row 23 is untouched by this edit)
$donorZero = $ws.Range("C23")   # shared text "0"
$donorNA = $ws.Range("E23")     # shared text "***.*"

# --- Update header title strings (volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  13"
$ws.Range("C9").Value = "Report Covering the Week  3/27/2023  Through  4/2/2023"

# --- Update weekly crime statistics table (rows 15-29) ---

# Row 15
$donorZero.Copy($ws.Range("C15"))
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -66.666666666666
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = -60
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -33.333333333333

# Row 16
$ws.Range("C16").Value = 3
$ws.Range("E16").Value = -50
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = -59.090909090909
$ws.Range("I16").Value = 46
$ws.Range("J16").Value = 60
$ws.Range("K16").Value = -23.333333333333
$ws.Range("L16").Value = 12.195121951219
$ws.Range("M16").Value = 24.324324324324
$ws.Range("N16").Value = -78.403755868544

# Row 17
$ws.Range("C17").Value = 2
$donorZero.Copy($ws.Range("D17"))
$donorNA.Copy($ws.Range("E17"))
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = 8.333333333333
$ws.Range("I17").Value = 40
$ws.Range("K17").Value = 25
$ws.Range("L17").Value = 53.846153846153
$ws.Range("M17").Value = 81.818181818181
$ws.Range("N17").Value = -18.367346938775

# Row 18
$ws.Range("C18").Value = 10
$ws.Range("D18").Value = 9
$ws.Range("E18").Value = 11.111111111111
$ws.Range("F18").Value = 24
$ws.Range("G18").Value = 37
$ws.Range("H18").Value = -35.135135135135
$ws.Range("I18").Value = 76
$ws.Range("J18").Value = 85
$ws.Range("K18").Value = -10.588235294117
$ws.Range("L18").Value = 38.181818181818
$ws.Range("M18").Value = 33.333333333333
$ws.Range("N18").Value = -58.695652173913

# Row 19
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = -10
$ws.Range("F19").Value = 84
$ws.Range("G19").Value = 85
$ws.Range("H19").Value = -1.176470588235
$ws.Range("I19").Value = 309
$ws.Range("J19").Value = 252
$ws.Range("K19").Value = 22.619047619047
$ws.Range("L19").Value = 145.238095238095
$ws.Range("M19").Value = 30.379746835443
$ws.Range("N19").Value = -44.524236983842

# Row 20
$donorZero.Copy($ws.Range("D20"))
$donorNA.Copy($ws.Range("E20"))
$ws.Range("I20").Value = 7
$ws.Range("K20").Value = -30
$ws.Range("L20").Value = 133.333333333333
$ws.Range("M20").Value = -12.5
$ws.Range("N20").Value = -96.551724137931

# Row 21
$ws.Range("C21").Value = 34
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = -5.555555555555
$ws.Range("F21").Value = 133
$ws.Range("G21").Value = 161
$ws.Range("H21").Value = -17.391304347826
$ws.Range("I21").Value = 480
$ws.Range("J21").Value = 444
$ws.Range("K21").Value = 8.108108108108
$ws.Range("L21").Value = 91.235059760956
$ws.Range("M21").Value = 32.231404958677
$ws.Range("N21").Value = -60.297766749379

# Row 22
$donorZero.Copy($ws.Range("D22"))
$donorNA.Copy($ws.Range("E22"))
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -60
$ws.Range("M22").Value = -38.095238095238

# Row 24
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 44
$ws.Range("E24").Value = -27.272727272727
$ws.Range("F24").Value = 119
$ws.Range("G24").Value = 132
$ws.Range("H24").Value = -9.848484848484
$ws.Range("I24").Value = 429
$ws.Range("J24").Value = 405
$ws.Range("K24").Value = 5.925925925925
$ws.Range("L24").Value = 57.720588235294
$ws.Range("M24").Value = 22.222222222222

# Row 25
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -42.857142857142
$ws.Range("F25").Value = 35
$ws.Range("G25").Value = 42
$ws.Range("H25").Value = -16.666666666666
$ws.Range("I25").Value = 102
$ws.Range("J25").Value = 76
$ws.Range("K25").Value = 34.210526315789
$ws.Range("L25").Value = 155
$ws.Range("M25").Value = 100

# Row 26
$donorZero.Copy($ws.Range("C26"))
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = -100
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -75
$ws.Range("J26").Value = 7
$ws.Range("K26").Value = -42.857142857142

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = -70
$ws.Range("I27").Value = 20
$ws.Range("J27").Value = 16
$ws.Range("K27").Value = 25
$ws.Range("L27").Value = 185.714285714286

# Row 28
$ws.Range("L28").Value = -100

# Row 29
$ws.Range("L29").Value = -100
